$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "61.406.19"
$ws.Range("E2").Value = "  +0.42%  "

$ws.Range("D3").Value = "3.381.09"
$ws.Range("E3").Value = "  +0.30%  "

$ws.Range("E4").Value = "  +0.08%  "

$ws.Range("D5").Value = "'575.65"
$ws.Range("E5").Value = "  +0.89%  "

$ws.Range("D6").Value = "'137.00"
$ws.Range("E6").Value = "  +0.67%  "

$ws.Range("E7").Value = "  +0.03%  "

$ws.Range("D8").Value = "3.380.83"
$ws.Range("E8").Value = "  +0.24%  "

$ws.Range("D9").Value = "'0.474"
$ws.Range("E9").Value = "  -0.27%  "

$ws.Range("D10").Value = "'7.48"
$ws.Range("E10").Value = "  -1.91%  "

$ws.Range("E11").Value = "  +1.88%  "

$ws.Range("E12").Value = "  -0.37%  "

$ws.Range("D13").Value = "3.956.19"
$ws.Range("E13").Value = "  +0.43%  "

$ws.Range("E14").Value = "  +2.45%  "

$ws.Range("E15").Value = "  +2.01%  "

$ws.Range("D16").Value = "'26.01"
$ws.Range("E16").Value = "  +3.30%  "

$ws.Range("D17").Value = "3.379.45"
$ws.Range("E17").Value = "  +0.31%  "

$ws.Range("D18").Value = "61.488.77"
$ws.Range("E18").Value = "  +0.60%  "

$ws.Range("D19").Value = "'14.09"
$ws.Range("E19").Value = "  +0.49%  "

$ws.Range("E20").Value = "  +1.05%  "

$ws.Range("E21").Value = "  -0.16%  "

$ws.Range("D22").Value = "'376.90"
$ws.Range("E22").Value = "  +0.22%  "

$ws.Range("E23").Value = "  -3.05%  "

$ws.Range("D24").Value = "3.517.60"
$ws.Range("E24").Value = "  +0.51%  "

$ws.Range("D25").Value = "'1.00"
$ws.Range("E25").Value = "  -0.10%  "

$ws.Range("E26").Value = "  +7.24%  "

$ws.Range("D27").Value = "'71.58"
$ws.Range("E27").Value = "  +1.11%  "

$ws.Range("E28").Value = "  +5.68%  "

$ws.Range("E29").Value = "  -3.50%  "

$ws.Range("E30").Value = "  +0.11%  "

$ws.Range("D31").Value = "'8.26"
$ws.Range("E31").Value = "  +1.86%  "

$ws.Range("E32").Value = "  +3.39%  "

$ws.Range("E33").Value = "  +1.88%  "

$ws.Range("E34").Value = "  +0.14%  "

$ws.Range("D35").Value = "'23.47"
$ws.Range("E35").Value = "  -0.18%  "

$ws.Range("D36").Value = "'5.29"
$ws.Range("E36").Value = "  -4.80%  "

$ws.Range("E37").Value = "  -1.87%  "

$ws.Range("E38").Value = "  -1.49%  "

$ws.Range("D39").Value = "'165.03"
$ws.Range("E39").Value = "  +0.58%  "

$ws.Range("E40").Value = "  -2.31%  "

$ws.Range("E41").Value = "  -0.09%  "

$ws.Range("B42").Value = "Stacks"
$ws.Range("C42").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D42").Value = "'1.72"
$ws.Range("E42").Value = "  +6.65%  "

$ws.Range("B43").Value = "Mantle"
$ws.Range("C43").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D43").Value = "'0.775"
$ws.Range("E43").Value = "  +2.01%  "

$ws.Range("E44").Value = "  +0.29%  "

$ws.Range("B45").Value = "OKB"
$ws.Range("C45").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D45").Value = "'41.50"
$ws.Range("E45").Value = "  -0.10%  "

$ws.Range("B46").Value = "Filecoin"
$ws.Range("C46").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D46").Value = "'4.41"
$ws.Range("E46").Value = "  +0.13%  "

$ws.Range("D47").Value = "'24.76"
$ws.Range("E47").Value = "  +6.94%  "

$ws.Range("E48").Value = "  -2.15%  "

$ws.Range("D49").Value = "'22.76"
$ws.Range("E49").Value = "  -1.46%  "

$ws.Range("D50").Value = "2.348.51"
$ws.Range("E50").Value = "  +4.22%  "

$ws.Range("E51").Value = "  +1.73%  "
